$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 17-19 need column-A header styling (bold, centered, bordered)
# to match the existing style used by rows 2-16 in column A (style index 1).
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("A16").Copy($ws.Range("A18"))
$ws.Range("A16").Copy($ws.Range("A19"))

# Row 10: Gaussian-Quadrature
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.4592399385408076
$ws.Range("D10").Value = 2.343758461795015
$ws.Range("E10").Value = 1.226217356742332
$ws.Range("F10").Value = 0.8846105184970876
$ws.Range("G10").Value = 0.4592399385408076
$ws.Range("H10").Value = 2.343758461795015
$ws.Range("I10").Value = 0.9295949883073978
$ws.Range("J10").Value = 1.062061652702853
$ws.Range("K10").Value = 0.6296575682396935
$ws.Range("L10").Value = 1.396241255662679
$ws.Range("M10").Value = 0.4592399385408076
$ws.Range("N10").Value = 1.784987909268674
$ws.Range("O10").Value = 1.22845656889381
$ws.Range("P10").Value = 1.116422717560983

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.08887436071823981
$ws.Range("D11").Value = 2.248172070589836
$ws.Range("E11").Value = 1.250525021303518
$ws.Range("F11").Value = 0.9736212635553039
$ws.Range("G11").Value = 0.08887436071823981
$ws.Range("H11").Value = 2.248172070589836
$ws.Range("I11").Value = 0.7582528072239293
$ws.Range("J11").Value = 1.369194225814396
$ws.Range("K11").Value = 0.4434056787342174
$ws.Range("L11").Value = 1.597875905015875
$ws.Range("M11").Value = 0.08887436071823981
$ws.Range("N11").Value = 1.749348545946677
$ws.Range("O11").Value = 1.140298179041724
$ws.Range("P11").Value = 1.091240166619414

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.08907486232426359
$ws.Range("D12").Value = 2.235208833344132
$ws.Range("E12").Value = 1.253331735869397
$ws.Range("F12").Value = 0.9738047673636576
$ws.Range("G12").Value = 0.08907486232426359
$ws.Range("H12").Value = 2.235208833344132
$ws.Range("I12").Value = 0.7600415591664825
$ws.Range("J12").Value = 1.370895911859556
$ws.Range("K12").Value = 0.4442225355853669
$ws.Range("L12").Value = 1.594353079807956
$ws.Range("M12").Value = 0.08907486232426359
$ws.Range("N12").Value = 1.744270284606765
$ws.Range("O12").Value = 1.137855049725363
$ws.Range("P12").Value = 1.090116660665101

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.08888207364973628
$ws.Range("D13").Value = 2.244750961627946
$ws.Range("E13").Value = 1.250634607553329
$ws.Range("F13").Value = 0.9739282173878996
$ws.Range("G13").Value = 0.08888207364973628
$ws.Range("H13").Value = 2.244750961627946
$ws.Range("I13").Value = 0.7584324945385025
$ws.Range("J13").Value = 1.369775965878517
$ws.Range("K13").Value = 0.4436605206040584
$ws.Range("L13").Value = 1.597582400438417
$ws.Range("M13").Value = 0.08888207364973628
$ws.Range("N13").Value = 1.747692784590637
$ws.Range("O13").Value = 1.139548965054728
$ws.Range("P13").Value = 1.090955905209801

# Row 14: NoRotation-tilt60deg
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.05600400000000012
$ws.Range("D14").Value = 4.425283999999993
$ws.Range("E14").Value = 0.7695360000000021
$ws.Range("F14").Value = 0.9609279999999999
$ws.Range("G14").Value = 0.05600400000000012
$ws.Range("H14").Value = 4.425283999999993
$ws.Range("I14").Value = 0.4796480000000019
$ws.Range("J14").Value = 1.076976
$ws.Range("K14").Value = 0.2972360000000004
$ws.Range("L14").Value = 2.089796000000003
$ws.Range("M14").Value = 0.05600400000000012
$ws.Range("N14").Value = 2.597409999999998
$ws.Range("O14").Value = 1.552937999999999
$ws.Range("P14").Value = 1.269426

# Row 15: Rotation-NoTilt
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 6.880000000000007
$ws.Range("E15").Value = 0.01
$ws.Range("F15").Value = 1.0471
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 6.880000000000007
$ws.Range("I15").Value = 0.01
$ws.Range("J15").Value = 0.7793124999999994
$ws.Range("K15").Value = 0.15
$ws.Range("L15").Value = 2.859575000000008
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 3.445000000000003
$ws.Range("O15").Value = 1.984275000000002
$ws.Range("P15").Value = 1.466998437500002

# Row 16: Rotation-60detTilt
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.4149965072383993
$ws.Range("D16").Value = 4.382975252479986
$ws.Range("E16").Value = 0.4217958011904027
$ws.Range("F16").Value = 1.024550997196797
$ws.Range("G16").Value = 0.4149965072383993
$ws.Range("H16").Value = 4.382975252479986
$ws.Range("I16").Value = 0.4245937657856008
$ws.Range("J16").Value = 0.8864219228160024
$ws.Range("K16").Value = 0.4937311351808002
$ws.Range("L16").Value = 2.090844984217598
$ws.Range("M16").Value = 0.4150318049279991
$ws.Range("N16").Value = 2.402385526835194
$ws.Range("O16").Value = 1.561079639526396
$ws.Range("P16").Value = 1.267488795763198

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9868060623377021
$ws.Range("D17").Value = 0.9930789243222364
$ws.Range("E17").Value = 0.9974135557467769
$ws.Range("F17").Value = 0.9986870260195182
$ws.Range("G17").Value = 0.9868060623377021
$ws.Range("H17").Value = 0.9930789243222364
$ws.Range("I17").Value = 0.9991679169844967
$ws.Range("J17").Value = 0.9897679080926864
$ws.Range("K17").Value = 0.9931799029441626
$ws.Range("L17").Value = 1.002912026506156
$ws.Range("M17").Value = 0.9868390662738894
$ws.Range("N17").Value = 0.9952462400345067
$ws.Range("O17").Value = 0.9939963921065584
$ws.Range("P17").Value = 0.995126665369217

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.06873374379983
$ws.Range("D18").Value = 0.9104053923201822
$ws.Range("E18").Value = 0.9826430690872173
$ws.Range("F18").Value = 0.9811588643937925
$ws.Range("G18").Value = 1.06873374379983
$ws.Range("H18").Value = 0.9104053923201822
$ws.Range("I18").Value = 1.016727508926961
$ws.Range("J18").Value = 0.9711677033595348
$ws.Range("K18").Value = 1.03989666415187
$ws.Range("L18").Value = 0.9215144027830765
$ws.Range("M18").Value = 1.06873374379983
$ws.Range("N18").Value = 0.9465242307036998
$ws.Range("O18").Value = 0.9857352674002555
$ws.Range("P18").Value = 0.9865309186028081

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.25131244311135
$ws.Range("D19").Value = 0.7297360846234779
$ws.Range("E19").Value = 1.132659122391168
$ws.Range("F19").Value = 0.8661021080975179
$ws.Range("G19").Value = 1.25131244311135
$ws.Range("H19").Value = 0.7297360846234779
$ws.Range("I19").Value = 1.158653932353331
$ws.Range("J19").Value = 0.9233265154687598
$ws.Range("K19").Value = 1.03462407566255
$ws.Range("L19").Value = 0.7571216481206009
$ws.Range("M19").Value = 1.25137879839358
$ws.Range("N19").Value = 0.9311976035073228
$ws.Range("O19").Value = 0.9949524395558785
$ws.Range("P19").Value = 0.9816919912285944

